# Update the public EPEX spot price workbook with the latest day's data.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": append a new day column (X) with the hourly prices for
# 07-jul.
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

$wsPrix.Range("X1").Value = "07-jul"
# Reuse the existing header style (bold, centered, bordered) instead of
# letting Excel mint a brand-new style entry for the new header cell.
$wsPrix.Range("W1").Copy()
$wsPrix.Range("X1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$prixValues = @(46.35, 35.82, 38.15, 25.44, 29.2, 35.85, 40.46, 38.13, 58.71, 57.5, 37.5, 18.01, 26.99, 8.13, 8.7, 5, 3.6, 4.31, 34.47, 54.37, 50.2, 54.66, 92.27, 78.98)

for ($i = 0; $i -lt $prixValues.Length; $i++) {
    $row = $i + 2
    $wsPrix.Cells.Item($row, 24).Value = $prixValues[$i]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append two new daily rows (2025-07-05 / 2025-07-06).
# Force the date column to Text first so the ISO-looking labels are kept as
# literal strings instead of being auto-converted to date serials.
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

$wsGaz.Range("A21").NumberFormat = "@"
$wsGaz.Range("A21").Value = "2025-07-05"
$wsGaz.Range("B21").Value = 32.5

$wsGaz.Range("A22").NumberFormat = "@"
$wsGaz.Range("A22").Value = "2025-07-06"
$wsGaz.Range("B22").Value = 32.5

# ---------------------------------------------------------------------------
# Sheet "CO2": append two new daily rows (2025-07-05 / 2025-07-06).
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A21").NumberFormat = "@"
$wsCo2.Range("A21").Value = "2025-07-05"
$wsCo2.Range("B21").Value = 70.92

$wsCo2.Range("A22").NumberFormat = "@"
$wsCo2.Range("A22").Value = "2025-07-06"
$wsCo2.Range("B22").Value = 70.92
